# Append a new "Week 8" section to the end of the progress log, mirroring the
# structure used by the existing weekly sections (a Heading1 paragraph followed
# by a bulleted ListParagraph outline at levels 1/2).
$d = $word.ActiveDocument

# Each entry: @("Heading1" -or- "ListParagraph", <list level: 0 = not a list item>, <text>)
$items = @(
    @("Heading1", 0, "Week 8 (25 hours)"),
    @("ListParagraph", 1, "Reworked host menu to be more intuitive"),
    @("ListParagraph", 2, "Both game modes shown at once"),
    @("ListParagraph", 2, "User-friendly text on the sliders"),
    @("ListParagraph", 2, "Made alignment and font consistent"),
    @("ListParagraph", 1, "Reworked combat messages"),
    @("ListParagraph", 2, "Removed trivial events from showing on the log (like knocking a limb off)"),
    @("ListParagraph", 2, "Added variations of each message to randomly choose from"),
    @("ListParagraph", 2, "Made a distinction between bombing a player and killing a player with a bomb"),
    @("ListParagraph", 2, "Added feedback messages on death telling the player how they died"),
    @("ListParagraph", 1, "Expanded the in-game menu with How to Play, Controls, and Options menus"),
    @("ListParagraph", 2, "Options can now be applied on the fly "),
    @("ListParagraph", 1, "Players are greeted with a game rules screen when they join a server, then click to spawn"),
    @("ListParagraph", 2, "The host can now also see the rules of each mode from the host screen"),
    @("ListParagraph", 1, "Various issues addressed in response to player feedback"),
    @("ListParagraph", 2, "Added a motion blur toggle to the options menu "),
    @("ListParagraph", 2, "Increased the blurred effect when losing your skull"),
    @("ListParagraph", 2, "Increased new round delay to give players more time"),
    @("ListParagraph", 2, "Bombs can now kill players if they have no arms or legs"),
    @("ListParagraph", 2, "Added sound effects to the bomb coffins – opening and closing, collecting bombs")
)

# Find the paragraph that ends the Week 7 section ("...last selected button"),
# which is the anchor after which the new content is inserted.
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains("last selected button")) {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph containing 'last selected button'"
}

$idx = $anchorIndex
foreach ($item in $items) {
    $style = $item[0]
    $level = $item[1]
    $text = $item[2]

    $d.Paragraphs.Item($idx).Range.InsertParagraphAfter() | Out-Null
    $idx = $idx + 1
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Range.Text = $text
    $newPara.Style = $style
    if ($level -gt 0) {
        $newPara.Range.ListFormat.ListLevelNumber = $level
    }
}

Write-Host "Inserted" ($idx - $anchorIndex) "paragraphs after index" $anchorIndex
